# Add 2022 data column (new column G holds 2022 counts, column F keeps 2021),
# plus a new "Total " column H that sums each row across years B:G.
# Also corrects a couple of previously-missing 2021 values and one 2020 value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / fill in 2021 (column F) values that changed or were missing ---
$ws.Range("F2").Value = 3
$ws.Range("F5").Value = 8
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 0

# --- New 2022 header + data (column G) ---
$ws.Range("G1").Value = 2022

$g2022 = @{
    2  = 1
    3  = 4
    4  = 4
    5  = 10
    6  = 12
    7  = 5
    8  = 3
    9  = 5
    10 = 2
    11 = 3
    12 = 4
    13 = 0
}
foreach ($row in $g2022.Keys) {
    $ws.Cells.Item($row, 7).Value = $g2022[$row]
}

# --- One corrected 2020 value (column E) ---
$ws.Range("E12").Value = 5

# --- New "Total " column (H), row-wise sums across B:G ---
$ws.Range("H1").Value = "Total "
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 8).Formula = "=SUM(B" + $row + ":G" + $row + ")"
}

# --- Row 14 grand totals: convert static numbers to SUM formulas, extend
#     through the new F/G/H columns ---
$ws.Range("B14").Formula = "=SUM(B2:B13)"
$ws.Range("C14").Formula = "=SUM(C2:C13)"
$ws.Range("D14").Formula = "=SUM(D2:D13)"
$ws.Range("E14").Formula = "=SUM(E2:E13)"
$ws.Range("F14").Formula = "=SUM(F2:F13)"
$ws.Range("G14").Formula = "=SUM(G2:G13)"
$ws.Range("H14").Formula = "=SUM(H2:H13)"

# --- Styling: center (horizontal + vertical) alignment on the new data
#     cells -- F2:F14 plus all of G1:H14. Build the style once on F2 and
#     fan it out with a format-only paste so every cell shares one style
#     entry instead of each assignment minting its own. ---
$ws.Range("F2").ClearFormats()
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4108
$ws.Range("F2").Copy()
$ws.Range("F3:F14").PasteSpecial(-4122)
$ws.Range("G1:H14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Final selection, matching the saved workbook state ---
$ws.Range("B15").Select()
